$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text representation
# (otherwise Excel would silently coerce numeric-looking strings to
# floating point numbers and lose formatting such as trailing zeros).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "26.536.76"
$ws.Range("E2").Value = "  +0.19%  "

$ws.Range("D3").Value = "1.815.88"
$ws.Range("E3").Value = "  +0.41%  "

$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.54%  "

$ws.Range("D5").Value = "1.003"
$ws.Range("E5").Value = "  -0.52%  "

$ws.Range("D6").Value = "306.36"

$ws.Range("D7").Value = "0.4533"
$ws.Range("E7").Value = "  -0.74%  "

$ws.Range("D8").Value = "0.3588"
$ws.Range("E8").Value = "  -2.10%  "

$ws.Range("D9").Value = "46.35"
$ws.Range("E9").Value = "  +2.80%  "

$ws.Range("D10").Value = "0.07088"
$ws.Range("E10").Value = "  -0.78%  "

$ws.Range("D11").Value = "0.8951"
$ws.Range("E11").Value = "  +2.06%  "

$ws.Range("D12").Value = "0.07776"
$ws.Range("E12").Value = "  -0.21%  "

$ws.Range("D13").Value = "19.32"
$ws.Range("E13").Value = "  -0.09%  "

$ws.Range("D14").Value = "1.837.69"
$ws.Range("E14").Value = "  +2.22%  "

$ws.Range("D15").Value = "5.265"
$ws.Range("E15").Value = "  -0.26%  "

$ws.Range("D16").Value = "6.315"
$ws.Range("E16").Value = "  -0.79%  "

$ws.Range("D17").Value = "85.29"
$ws.Range("E17").Value = "  -0.81%  "

$ws.Range("D18").Value = "1.005"
$ws.Range("E18").Value = "  -0.53%  "

$ws.Range("D19").Value = "0.000008603"
$ws.Range("E19").Value = "  +0.14%  "

$ws.Range("E20").Value = "  -0.53%  "

$ws.Range("D21").Value = "26.586.19"
$ws.Range("E21").Value = "  +0.25%  "

$ws.Range("D22").Value = "14.16"
$ws.Range("E22").Value = "  -0.69%  "

$ws.Range("D23").Value = "4.962"
$ws.Range("E23").Value = "  -0.78%  "

$ws.Range("D24").Value = "10.51"
$ws.Range("E24").Value = "  +0.58%  "

$ws.Range("D25").Value = "1.966"
$ws.Range("E25").Value = "  -0.89%  "

$ws.Range("D26").Value = "151.28"
$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("D27").Value = "17.78"
$ws.Range("E27").Value = "  -1.15%  "

$ws.Range("D28").Value = "2.047"
$ws.Range("E28").Value = "  -0.27%  "

$ws.Range("D29").Value = "112.24"
$ws.Range("E29").Value = "  -0.26%  "

$ws.Range("D30").Value = "4.844"
$ws.Range("E30").Value = "  -0.02%  "

$ws.Range("D31").Value = "0.08718"
$ws.Range("E31").Value = "  +0.49%  "

$ws.Range("D32").Value = "3.124"
$ws.Range("E32").Value = "  +2.16%  "

$ws.Range("D33").Value = "0.7573"
$ws.Range("E33").Value = "  +3.76%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "4.430"
$ws.Range("E34").Value = "  -0.80%  "

$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").Value = "2.721"
$ws.Range("E35").Value = "  +9.11%  "

$ws.Range("D36").Value = "1.113"
$ws.Range("E36").Value = "  -0.30%  "

$ws.Range("D37").Value = "1.074"

$ws.Range("D38").Value = "0.01937"
$ws.Range("E38").Value = "  +0.13%  "

$ws.Range("D39").Value = "2.917"
$ws.Range("E39").Value = "  +0.78%  "

$ws.Range("D40").Value = "0.05105"
$ws.Range("E40").Value = "  -0.10%  "

$ws.Range("D41").Value = "0.5112"
$ws.Range("E41").Value = "  +2.27%  "

$ws.Range("D42").Value = "6.770"
$ws.Range("E42").Value = "  -2.90%  "

$ws.Range("D43").Value = "0.1513"
$ws.Range("E43").Value = "  -3.14%  "

$ws.Range("D44").Value = "8.043"
$ws.Range("E44").Value = "  -0.96%  "

$ws.Range("D45").Value = "0.4705"
$ws.Range("E45").Value = "  +2.05%  "

$ws.Range("E46").Value = "  -0.65%  "

$ws.Range("D47").Value = "10.05"
$ws.Range("E47").Value = "  +1.02%  "

$ws.Range("D48").Value = "100.51"
$ws.Range("E48").Value = "  -0.43%  "

$ws.Range("D49").Value = "1.573"
$ws.Range("E49").Value = "  -0.98%  "

$ws.Range("D50").Value = "0.05988"
$ws.Range("E50").Value = "  -0.48%  "

$ws.Range("D51").Value = "63.90"
$ws.Range("E51").Value = "  -0.27%  "
